$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1) Add the new "Sheet2" right after "Sheet1".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# 2) Populate Sheet2 with the rows that used to live on Sheet1 (old rows 3-5),
#    giving the ProjectName column (A) its new values. Introduce the brand new
#    shared strings in the same order they appear in the target workbook so the
#    shared-strings table compacts into the expected order.
$sheet2Rows = @(
    @("Online Test2", -1, 3, "Laura's Test Environment", "Convert to UpdatePackage", "Martha", "Same", "DB"),
    @("OnlineTest 3", 1, 3, "Laura's Build Environment", "Update Contact page", "Bocephus", "Tooty", "Back End"),
    @("Online Test ", 0, 3, "BR549", "Test Register function", "Marty", "Sarah", "Architecture")
)

for ($r = 0; $r -lt $sheet2Rows.Length; $r++) {
    $rowValues = $sheet2Rows[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws2.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

# 3) Rename the project on row 2 of Sheet1.
$ws1.Range("A2").Value = "Tuesday 24th"

# 4) Drop the now-duplicated rows 3-5 from Sheet1 (their data now lives on Sheet2).
$ws1.Rows("3:5").Delete()

# 5) Fix up the selections so Sheet1 stays the active tab, and both sheets show
#    the selection the original workbook ended up with.
$ws2.Range("A1:H3").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A2").Select() | Out-Null
